# allocation rule updated with 5 and 10 mi rad
#
# Columns E (pop_sq_mile_3mi numerator helper) and G on rows 2-10 switch from
# numeric storage to text storage of the very same number (the author started
# formatting these as text), and columns R (total_risk) / S (total_risk_resp)
# get refreshed values reflecting the new 5/10 mi radius allocation rule.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $text) {
    # Force the cell to be written back as a shared-string (text) cell
    # holding the exact same characters it already displayed, then drop the
    # "@" text number-format again so no extra style is left behind on the
    # cell (matches a plain <c t="s"> with no s="" attribute).
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).ClearFormats()
}

Set-TextValue $ws "E2" "79.1835308238153"
Set-TextValue $ws "G2" "172.624280046488"
$ws.Range("R2").Value = 43
$ws.Range("S2").Value = 0.53

Set-TextValue $ws "E3" "103.612502034919"
Set-TextValue $ws "G3" "45.6605130373706"
$ws.Range("R3").Value = 50
$ws.Range("S3").Value = 0.6

Set-TextValue $ws "E4" "60.2481873222828"
Set-TextValue $ws "G4" "1262.07946461947"
$ws.Range("R4").Value = 30
$ws.Range("S4").Value = 0.404347826086957

Set-TextValue $ws "E5" "50.1163080888244"
Set-TextValue $ws "G5" "1478.2014642559"
$ws.Range("R5").Value = 30
$ws.Range("S5").Value = 0.384615384615385

Set-TextValue $ws "E6" "155.013984349669"
Set-TextValue $ws "G6" "58.0399248348152"
$ws.Range("R6").Value = 34.2857142857143
$ws.Range("S6").Value = 1.37142857142857

Set-TextValue $ws "E7" "163.018378976147"
Set-TextValue $ws "G7" "151.381704044615"
$ws.Range("R7").Value = 128.75
$ws.Range("S7").Value = 0.5375

Set-TextValue $ws "E8" "194.534577199069"
Set-TextValue $ws "G8" "124.414900160566"
$ws.Range("R8").Value = 127
$ws.Range("S8").Value = 0.57

Set-TextValue $ws "E9" "41.168584416499"
Set-TextValue $ws "G9" "828.37436563237"
$ws.Range("R9").Value = 21.3333333333333
$ws.Range("S9").Value = 0.266666666666667

Set-TextValue $ws "E10" "317.643922321843"
Set-TextValue $ws "G10" "57.7501998650348"
$ws.Range("R10").Value = 20
$ws.Range("S10").Value = 0.2
